$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ADP")

$ws.Range("D45").Value = 27665800
$ws.Range("D46").Value = 31823300
$ws.Range("D48").Value = 1587400
$ws.Range("D52").Value = 3069500
$ws.Range("D54").Value = 38849100
$ws.Range("D59").Value = 30277300
$ws.Range("D60").Value = 30412700
$ws.Range("D62").Value = 1698100
$ws.Range("D66").Value = 34113200
$ws.Range("D72").Value = 16546600
$ws.Range("D76").Value = 4735900

$ws.Range("D91").Value = -206100
$ws.Range("E91").Value = -240200
$ws.Range("F91").Value = -168500
$ws.Range("G91").Value = -158800
$ws.Range("H91").Value = -159800
$ws.Range("I91").Value = -130300
$ws.Range("J91").Value = -140100
